# "integrated gui and api test"
# The "Global" sheet stores global test arguments. Add the
# "arg_ReservationNumber" argument name in A1 and best-fit column A to it,
# matching how SpreadsheetGear/Excel auto-sizes a column after typing a
# long value into an otherwise-default-width column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Global")

$ws.Range("A1").Value = "arg_ReservationNumber"
$ws.Columns.Item(1).AutoFit()
